# Auto-generated edit script: updates computed price/profit columns (H-N)
# on multiple rows across multiple sheets, per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 473.42105
$ws.Range("I9").Value = 408.2143
$ws.Range("K9").Value = 408.2143
$ws.Range("M9").Value = -239.2143

$ws.Range("H32").Value = 3211.0588
$ws.Range("I32").Value = 2550
$ws.Range("J32").Value = 4797.6
$ws.Range("K32").Value = 2550
$ws.Range("L32").Value = 4797.6
$ws.Range("M32").Value = -2224
$ws.Range("N32").Value = -5449.6

$ws.Range("H98").Value = 2707.7297
$ws.Range("I98").Value = 2707.7297
$ws.Range("K98").Value = 2707.7297
$ws.Range("M98").Value = -1209.7297

$ws.Range("H113").Value = 80568860
$ws.Range("J113").Value = 107160860
$ws.Range("L113").Value = 107160860
$ws.Range("N113").Value = -107167368

$ws.Range("H122").Value = 2707.7297
$ws.Range("I122").Value = 2707.7297
$ws.Range("K122").Value = 8123.1891
$ws.Range("M122").Value = -5673.1891

$ws.Range("H129").Value = 1213.8334
$ws.Range("J129").Value = 2072.3333
$ws.Range("L129").Value = 6216.999899999999
$ws.Range("N129").Value = -16216.9999

$ws.Range("H138").Value = 2243.87
$ws.Range("I138").Value = 1635.15
$ws.Range("J138").Value = 2649.6833
$ws.Range("K138").Value = 4905.450000000001
$ws.Range("L138").Value = 7949.0499
$ws.Range("M138").Value = 234.5499999999993
$ws.Range("N138").Value = -18229.0499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 23812178
$ws.Range("I2").Value = 1209.25
$ws.Range("J2").Value = 100007280
$ws.Range("K2").Value = 1209.25
$ws.Range("L2").Value = 100007280
$ws.Range("M2").Value = -1096.25
$ws.Range("N2").Value = -100007506

$ws.Range("H7").Value = 60000
$ws.Range("J7").Value = 60000
$ws.Range("L7").Value = 60000
$ws.Range("N7").Value = -60228

$ws.Range("H61").Value = 28574206
$ws.Range("I61").Value = 2010.037
$ws.Range("K61").Value = 2010.037
$ws.Range("M61").Value = -1798.037

$ws.Range("H102").Value = 5885347.5
$ws.Range("I102").Value = 6454155.5
$ws.Range("K102").Value = 6454155.5
$ws.Range("M102").Value = -6452533.5

$ws.Range("H116").Value = 23812178
$ws.Range("I116").Value = 1209.25
$ws.Range("J116").Value = 100007280
$ws.Range("K116").Value = 1209.25
$ws.Range("L116").Value = 100007280
$ws.Range("M116").Value = 1084.75
$ws.Range("N116").Value = -100011868

$ws.Range("H122").Value = 2990.2
$ws.Range("I122").Value = 2237.0476
$ws.Range("K122").Value = 6711.1428
$ws.Range("M122").Value = -4261.1428

$ws.Range("H132").Value = 5096.212
$ws.Range("I132").Value = 2926.5
$ws.Range("K132").Value = 8779.5
$ws.Range("M132").Value = -6249.5

$ws.Range("H136").Value = 28574206
$ws.Range("I136").Value = 2010.037
$ws.Range("K136").Value = 6030.111
$ws.Range("M136").Value = -3480.111

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 64893.75
$ws.Range("J2").Value = 64893.75
$ws.Range("L2").Value = 64893.75
$ws.Range("N2").Value = -65119.75

$ws.Range("H3").Value = 23812178
$ws.Range("I3").Value = 1209.25
$ws.Range("J3").Value = 100007280
$ws.Range("K3").Value = 1209.25
$ws.Range("L3").Value = 100007280
$ws.Range("M3").Value = -1095.25
$ws.Range("N3").Value = -100007508

$ws.Range("H13").Value = 78499.664
$ws.Range("J13").Value = 78499.664
$ws.Range("L13").Value = 78499.664
$ws.Range("N13").Value = -78835.664

$ws.Range("H48").Value = 250000
$ws.Range("J48").Value = 250000
$ws.Range("L48").Value = 250000
$ws.Range("N48").Value = -250830

$ws.Range("H99").Value = 11114389
$ws.Range("I99").Value = 3250.25
$ws.Range("K99").Value = 3250.25
$ws.Range("M99").Value = -1752.25

$ws.Range("H118").Value = 60403.75
$ws.Range("J118").Value = 60403.75
$ws.Range("L118").Value = 60403.75
$ws.Range("N118").Value = -63717.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 235.1579
$ws.Range("I7").Value = 278.2
$ws.Range("J7").Value = 73.75
$ws.Range("K7").Value = 278.2
$ws.Range("L7").Value = 73.75
$ws.Range("M7").Value = -165.2
$ws.Range("N7").Value = -299.75

$ws.Range("H22").Value = 699.8182
$ws.Range("I22").Value = 824.75
$ws.Range("J22").Value = 366.66666
$ws.Range("K22").Value = 824.75
$ws.Range("L22").Value = 366.66666
$ws.Range("M22").Value = -474.75
$ws.Range("N22").Value = -1066.66666

$ws.Range("H31").Value = 5495.405
$ws.Range("I31").Value = 1203.9166
$ws.Range("K31").Value = 1203.9166
$ws.Range("M31").Value = -908.9166

$ws.Range("H34").Value = 5495.405
$ws.Range("I34").Value = 1203.9166
$ws.Range("K34").Value = 1203.9166
$ws.Range("M34").Value = -1001.9166

$ws.Range("H122").Value = 3809.3076
$ws.Range("I122").Value = 2780.4
$ws.Range("K122").Value = 8341.200000000001
$ws.Range("M122").Value = -5891.200000000001

$ws.Range("H132").Value = 3876.907
$ws.Range("I132").Value = 2190.3333
$ws.Range("K132").Value = 6570.999899999999
$ws.Range("M132").Value = -4040.999899999999

$ws.Range("H134").Value = 2802.513
$ws.Range("I134").Value = 1416.7667
$ws.Range("K134").Value = 4250.300099999999
$ws.Range("M134").Value = -1715.300099999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 52186.7
$ws.Range("J131").Value = 54854.42
$ws.Range("L131").Value = 164563.26
$ws.Range("N131").Value = -174643.26

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 55000
$ws.Range("J51").Value = 55000
$ws.Range("L51").Value = 55000
$ws.Range("N51").Value = -56018

$ws.Range("H118").Value = 27000
$ws.Range("J118").Value = 27000
$ws.Range("L118").Value = 27000
$ws.Range("N118").Value = -30314

$ws.Range("H126").Value = 4777
$ws.Range("I126").Value = 2432.1
$ws.Range("K126").Value = 7296.299999999999
$ws.Range("M126").Value = -4826.299999999999

$ws.Range("H132").Value = 2777
$ws.Range("I132").Value = 2078.0833
$ws.Range("J132").Value = 4873.75
$ws.Range("K132").Value = 6234.249899999999
$ws.Range("L132").Value = 14621.25
$ws.Range("M132").Value = -3704.249899999999
$ws.Range("N132").Value = -19681.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 861.04
$ws.Range("J16").Value = 838.8
$ws.Range("L16").Value = 838.8
$ws.Range("N16").Value = -1178.8

$ws.Range("H40").Value = 5861.222
$ws.Range("I40").Value = 4375.25
$ws.Range("K40").Value = 4375.25
$ws.Range("M40").Value = -4239.25

$ws.Range("H61").Value = 3452244
$ws.Range("I61").Value = 5558077
$ws.Range("J61").Value = 6335.364
$ws.Range("K61").Value = 5558077
$ws.Range("L61").Value = 6335.364
$ws.Range("M61").Value = -5557875
$ws.Range("N61").Value = -6739.364

$ws.Range("H113").Value = 3452244
$ws.Range("I113").Value = 5558077
$ws.Range("J113").Value = 6335.364
$ws.Range("K113").Value = 5558077
$ws.Range("L113").Value = 6335.364
$ws.Range("M113").Value = -5555907
$ws.Range("N113").Value = -10675.364

$ws.Range("H122").Value = 2735.2046
$ws.Range("I122").Value = 2081.6667
$ws.Range("K122").Value = 6245.000100000001
$ws.Range("M122").Value = -3795.000100000001

$ws.Range("H132").Value = 7940489
$ws.Range("I132").Value = 13160682
$ws.Range("K132").Value = 39482046
$ws.Range("M132").Value = -39479516

$ws.Range("H136").Value = 8244.593000000001
$ws.Range("I136").Value = 3318.4333
$ws.Range("K136").Value = 9955.2999
$ws.Range("M136").Value = -7405.2999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H60").Value = 66666
$ws.Range("J60").Value = 66666
$ws.Range("L60").Value = 66666
$ws.Range("N60").Value = -68310

$ws.Range("H100").Value = 1759.6
$ws.Range("I100").Value = 932.6667
$ws.Range("K100").Value = 1865.3334
$ws.Range("M100").Value = -1324.3334

$ws.Range("H117").Value = 56117
$ws.Range("J117").Value = 56117
$ws.Range("L117").Value = 56117
$ws.Range("N117").Value = -65295

$ws.Range("H132").Value = 4911.7
$ws.Range("I132").Value = 4586.7144
$ws.Range("K132").Value = 13760.1432
$ws.Range("M132").Value = -11230.1432
